$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2742.7144
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 3066.5
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 3066.5
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -3718.5

$ws.Range("H106").Value = 2689.9
$ws.Range("I106").Value = 2916.5
$ws.Range("J106").Value = 2350
$ws.Range("K106").Value = 2916.5
$ws.Range("L106").Value = 2350
$ws.Range("M106").Value = -2285.5
$ws.Range("N106").Value = -3612

$ws.Range("H132").Value = 2136.5273
$ws.Range("I132").Value = 1979.3541
$ws.Range("J132").Value = 3214.2856
$ws.Range("K132").Value = 5938.0623
$ws.Range("L132").Value = 9642.856800000001
$ws.Range("M132").Value = -3408.0623
$ws.Range("N132").Value = -14702.8568

$ws.Range("H139").Value = 52570
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 52570
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 52570
$ws.Range("N139").Value = -62850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10419.019
$ws.Range("I32").Value = 4835.7446
$ws.Range("J32").Value = 54154.668
$ws.Range("K32").Value = 4835.7446
$ws.Range("L32").Value = 54154.668
$ws.Range("M32").Value = -4548.7446
$ws.Range("N32").Value = -54728.668

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10900

$ws.Range("H139").Value = 45710
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45710
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 45710
$ws.Range("N139").Value = -55990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 137372
$ws.Range("I134").Value = 4066.6667
$ws.Range("J134").Value = 337330
$ws.Range("K134").Value = 12200.0001
$ws.Range("L134").Value = 1011990
$ws.Range("M134").Value = -9665.000100000001
$ws.Range("N134").Value = -1017060

$ws.Range("H138").Value = 54680
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 54680
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 54680
$ws.Range("N138").Value = -64960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1251.2632
$ws.Range("I5").Value = 955.8484999999999
$ws.Range("J5").Value = 3201
$ws.Range("K5").Value = 2867.5455
$ws.Range("L5").Value = 9603
$ws.Range("M5").Value = -2755.5455
$ws.Range("N5").Value = -9827

$ws.Range("H76").Value = 4333.3335
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 15000
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -15766

$ws.Range("H79").Value = 4333.3335
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -17652

$ws.Range("H80").Value = 5358.3335
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5358.3335
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 16075.0005
$ws.Range("N80").Value = -17947.0005

$ws.Range("H83").Value = 5358.3335
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5358.3335
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 48225.0015
$ws.Range("N83").Value = -57585.0015

$ws.Range("H87").Value = 14645.363
$ws.Range("I87").Value = 2950
$ws.Range("J87").Value = 21328.428
$ws.Range("K87").Value = 8850
$ws.Range("L87").Value = 63985.284
$ws.Range("M87").Value = -7602
$ws.Range("N87").Value = -66481.284

$ws.Range("H90").Value = 14645.363
$ws.Range("I90").Value = 2950
$ws.Range("J90").Value = 21328.428
$ws.Range("K90").Value = 26550
$ws.Range("L90").Value = 191955.852
$ws.Range("M90").Value = -20310
$ws.Range("N90").Value = -204435.852

$ws.Range("H113").Value = 796.34485
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 796.34485
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2389.03455
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6729.03455

$ws.Range("H117").Value = 2141
$ws.Range("I117").Value = 677
$ws.Range("J117").Value = 2407.182
$ws.Range("K117").Value = 2031
$ws.Range("L117").Value = 7221.545999999999
$ws.Range("M117").Value = 1411
$ws.Range("N117").Value = -14105.546

$ws.Range("H120").Value = 11260.692
$ws.Range("I120").Value = 3374.75
$ws.Range("J120").Value = 14765.556
$ws.Range("K120").Value = 10124.25
$ws.Range("L120").Value = 44296.66800000001
$ws.Range("M120").Value = -5286.25
$ws.Range("N120").Value = -53972.66800000001

$ws.Range("H135").Value = 1251.2632
$ws.Range("I135").Value = 955.8484999999999
$ws.Range("J135").Value = 3201
$ws.Range("K135").Value = 8602.636499999999
$ws.Range("L135").Value = 28809
$ws.Range("M135").Value = -6067.636499999999
$ws.Range("N135").Value = -33879

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5314.2856
$ws.Range("I70").Value = 4700
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 4700
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4430
$ws.Range("N70").Value = -9540

$ws.Range("H73").Value = 5314.2856
$ws.Range("I73").Value = 4700
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 4700
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3764
$ws.Range("N73").Value = -10872

$ws.Range("H122").Value = 3001.25
$ws.Range("I122").Value = 2935.6667
$ws.Range("J122").Value = 3198
$ws.Range("K122").Value = 8807.000100000001
$ws.Range("L122").Value = 9594
$ws.Range("M122").Value = -6357.000100000001
$ws.Range("N122").Value = -14494

$ws.Range("H132").Value = 4563.8
$ws.Range("I132").Value = 6881
$ws.Range("J132").Value = 3019
$ws.Range("K132").Value = 20643
$ws.Range("L132").Value = 9057
$ws.Range("M132").Value = -18113
$ws.Range("N132").Value = -14117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 9950
$ws.Range("I34").Value = 4000
$ws.Range("J34").Value = 12925
$ws.Range("K34").Value = 4000
$ws.Range("L34").Value = 12925
$ws.Range("M34").Value = -3828
$ws.Range("N34").Value = -13269

$ws.Range("H61").Value = 2734.5557
$ws.Range("I61").Value = 2658
$ws.Range("J61").Value = 3002.5
$ws.Range("K61").Value = 2658
$ws.Range("L61").Value = 3002.5
$ws.Range("M61").Value = -2456
$ws.Range("N61").Value = -3406.5

$ws.Range("H68").Value = 3420
$ws.Range("I68").Value = 3725
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 3725
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -2976
$ws.Range("N68").Value = -3698

$ws.Range("H71").Value = 3420
$ws.Range("I71").Value = 3725
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 18625
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -14881
$ws.Range("N71").Value = -18488

$ws.Range("H113").Value = 2734.5557
$ws.Range("I113").Value = 2658
$ws.Range("J113").Value = 3002.5
$ws.Range("K113").Value = 2658
$ws.Range("L113").Value = 3002.5
$ws.Range("M113").Value = -488
$ws.Range("N113").Value = -7342.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2268
$ws.Range("I122").Value = 2268
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6804
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -4354

$ws.Range("H132").Value = 1585.8
$ws.Range("I132").Value = 1529.8462
$ws.Range("J132").Value = 1949.5
$ws.Range("K132").Value = 4589.5386
$ws.Range("L132").Value = 5848.5
$ws.Range("M132").Value = -2059.5386
$ws.Range("N132").Value = -10908.5

$ws.Range("H133").Value = 69785.836
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 69785.836
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 69785.836
$ws.Range("N133").Value = -79905.836
